# Fix the "5a 14" typo -> "5 a 14" in the age-range header cells (F1 and N1)
# on every sheet that has that age-bracket table (all sheets except the
# "Global semana epidemiologica" summary sheet).

$wb = $excel.ActiveWorkbook

$sheetNames = @("VRS", "Ad", "Parainfluenza", "Inf A", "Inf B", "Metapnemovirus")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F1").Value = "5 a 14"
    $ws.Range("N1").Value = "5 a 14"
}
